$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date (column C) for rows 2-18 from 2023-10-05 (45204)
# to 2023-10-08 (45207), keeping the existing date formatting/style intact.
$ws.Range("C2:C18").Value = 45207
